# Saldo_guide.xlsx update:
#  - Refresh the "Dt. Referencia" column (G) from 2024-05-29 (serial 45441)
#    to 2024-05-31 (serial 45443) for every data row.
#  - Update the one client whose balance changed (row 42, "Saldo Previsto"
#    / "Vl. Total" columns D & H) from 317.97 to 334.63.
#  - Rename the sheet to match the new export timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 257 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq 45441) {
        $cell.Value2 = 45443
    }
}

# Corrected balance for row 42 (Saldo Previsto / Vl. Total)
$ws.Cells.Item(42, 4).Value2 = 334.63
$ws.Cells.Item(42, 8).Value2 = 334.63

# New export file name baked into the sheet title
$ws.Name = "IClientBalance-20240531-090045-"
